$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabla_MS")
$ws.Activate()

# Add new row 18 with the extra variable "Enfermedad"
$ws.Range("A18").Value = "Enfermedad"
$ws.Range("B18").Value = "categórica"
$ws.Range("C18").Value = "Enfermedad principal en la que se trata con dicho medicamento"

# Match style (center/center) used by the other rows in this table (s="3")
$ws.Range("A18:C18").HorizontalAlignment = -4108
$ws.Range("A18:C18").VerticalAlignment = -4108

# Update the view to match: scrolled so row 12 / column C is the top-left
# visible cell, with C18 as the active selection.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C18").Select()
